{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Summary of changes being applied (see task diff):\n//   1. A leading space is typed at the very start of the document, and the\n//      \"_GoBack\" bookmark (which Word maintains to mark the last edit\n//      location) moves there - i.e. right before \"Programming\" in the\n//      title \"Programming For Design\".\n//   2. The \"_GoBack\" bookmark is removed from its old location (an empty\n//      paragraph right after \"Draw Eyes\").\n//   3. A new paragraph is added right after the \"Flow Chart \" paragraph,\n//      explaining what the flow chart shows.\n//   4. A new paragraph is added right before the \"DO\" paragraph (after\n//      \"SUEDO Code\"), explaining the pseudocode basics.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// --- 1 & 2. Move the \"_GoBack\" bookmark to the start of the document ---\n\n// The very first paragraph holds the title \"Programming For Design\".\nconst titleParagraph = paragraphs.items[0];\n\n// Type a leading space at the start of the title paragraph.\nconst docStart = titleParagraph.getRange(\"Start\");\ndocStart.insertText(\" \", \"Before\");\nawait context.sync();\n\n// Remove the bookmark from its previous location (if present).\nconst goBackExists = context.document.bookmarks.exists(\"_GoBack\");\nawait context.sync();\nif (goBackExists.value) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// Re-find the word \"Programming\" (now preceded by the inserted space) and\n// drop a collapsed \"_GoBack\" bookmark immediately before it.\nconst progResults = titleParagraph.search(\"Programming\", { matchCase: true });\nprogResults.load(\"items\");\nawait context.sync();\nconst beforeProgramming = progResults.items[0].getRange(\"Start\");\nbeforeProgramming.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 3. Insert explanatory paragraph after \"Flow Chart \" ---\n\nlet flowChartParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Flow Chart \") {\n    flowChartParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (flowChartParagraph) {\n  flowChartParagraph.insertParagraph(\n    \"This flow chart explains how the code works, I needed to draw everything each frame, and checking to see where each eye was, then moving them if they had reached a certain position\",\n    \"After\"\n  );\n  await context.sync();\n}\n\n// --- 4. Insert explanatory paragraph before \"DO\" (after \"SUEDO Code\") ---\n\nlet suedoParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"SUEDO Code\") {\n    suedoParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (suedoParagraph) {\n  suedoParagraph.insertParagraph(\n    \"The Basics for setting up how my code will work, drawing everything and checking where the eyes are if they have reached there min or max position, reverse the rotation until they hit their min or max position, I would also add a small amount to avoid them not moving due to them already hitting their min or max position\",\n    \"After\"\n  );\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n#\n# Summary of changes being applied (see task diff):\n#   1. A leading space is typed at the very start of the document, and the\n#      \"_GoBack\" bookmark (which Word maintains to mark the last edit\n#      location) moves there - i.e. right before \"Programming\" in the\n#      title \"Programming For Design\".\n#   2. The \"_GoBack\" bookmark is removed from its old location (an empty\n#      paragraph right after \"Draw Eyes\").\n#   3. A new paragraph is added right after the \"Flow Chart \" paragraph,\n#      explaining what the flow chart shows.\n#   4. A new paragraph is added right before the \"DO\" paragraph (after\n#      \"SUEDO Code\"), explaining the pseudocode basics.\n\n$d = $word.ActiveDocument\n\n# --- 1 & 2. Move the \"_GoBack\" bookmark to the start of the document ---\n\n# Type a leading space at the very start of the document.\n$docStart = $d.Range(0, 0)\n$docStart.InsertBefore(\" \")\n\n# Remove the bookmark from its previous location (if present).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Find \"Programming\" (now preceded by the inserted leading space) and drop a\n# collapsed \"_GoBack\" bookmark immediately before it.\n$progRange = $d.Content\n$progRange.Find.Execute(\"Programming\")\n$bmRange = $d.Range($progRange.Start, $progRange.Start)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n# --- 3. Insert explanatory paragraph after \"Flow Chart \" ---\n\n$flowChartRange = $d.Content\nif ($flowChartRange.Find.Execute(\"Flow Chart \")) {\n    $flowChartPara = $flowChartRange.Paragraphs(1)\n    $flowChartParaRange = $flowChartPara.Range\n    $insertPos = $flowChartParaRange.End\n    $flowChartParaRange.InsertParagraphAfter()\n    $newRange = $d.Range($insertPos, $insertPos)\n    $newRange.InsertAfter(\"This flow chart explains how the code works, I needed to draw everything each frame, and checking to see where each eye was, then moving them if they had reached a certain position\")\n}\n\n# --- 4. Insert explanatory paragraph before \"DO\" (after \"SUEDO Code\") ---\n\n$suedoRange = $d.Content\nif ($suedoRange.Find.Execute(\"SUEDO Code\")) {\n    $suedoPara = $suedoRange.Paragraphs(1)\n    $suedoParaRange = $suedoPara.Range\n    $insertPos2 = $suedoParaRange.End\n    $suedoParaRange.InsertParagraphAfter()\n    $newRange2 = $d.Range($insertPos2, $insertPos2)\n    $newRange2.InsertAfter(\"The Basics for setting up how my code will work, drawing everything and checking where the eyes are if they have reached there min or max position, reverse the rotation until they hit their min or max position, I would also add a small amount to avoid them not moving due to them already hitting their min or max position\")\n}\n"}
